$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1 and Q1, matching the style used by the other header cells
# (bold font, thin border all around, centered horizontally, top vertically - same as O1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$headerRange = $ws.Range("P1:Q1")
$headerRange.Borders.LineStyle = 1         # xlContinuous (thin border)
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop

# Update existing columns I, K, M, O for data rows 2..25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # column I
    $ws.Cells.Item($r, 11).Value = 1  # column K
    $ws.Cells.Item($r, 13).Value = 2  # column M
    $ws.Cells.Item($r, 15).Value = 1  # column O

    # Add new columns P and Q
    $ws.Cells.Item($r, 16).Value = 2  # column P
    $ws.Cells.Item($r, 17).Value = 2  # column Q
}
